{"js": "// Remove the empty paragraph, the page-break paragraph, and the\n// \"\u00a9 2020 ...\" footer paragraph that immediately follow the final\n// bibliography entry (\"Edgard Blucher Ltda, 1993...\").\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst marker = \"Edgard Blucher Ltda, 1993.SILVA, R. R.; BOCCHI, N. ; ROCHA FILHO, R. C. Introdu\u00e7\u00e3o a qu\u00edmica experimental. S\u00e3o Paulo: McGraw-Hill, 1990.\";\nconst footer = \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\";\n\nconst items = paragraphs.items;\nlet markerIndex = -1;\nlet footerIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === marker) {\n    markerIndex = i;\n  }\n  if (items[i].text === footer) {\n    footerIndex = i;\n    break;\n  }\n}\n\nif (markerIndex === -1 || footerIndex === -1 || footerIndex <= markerIndex) {\n  throw new Error(\"Could not locate the paragraphs to delete.\");\n}\n\n// Delete every paragraph strictly after the bibliography marker up to and\n// including the footer paragraph (the blank paragraph, the page-break\n// paragraph, and the footer paragraph itself).\nfor (let i = footerIndex; i > markerIndex; i--) {\n  items[i].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the empty paragraph, the page-break paragraph, and the\n# \"(c) 2020 ...\" footer paragraph that immediately follow the final\n# bibliography entry (\"Edgard Blucher Ltda, 1993...\").\n$d = $word.ActiveDocument\n\n$marker = \"Edgard Blucher Ltda, 1993.SILVA, R. R.; BOCCHI, N. ; ROCHA FILHO, R. C. Introdu\u00e7\u00e3o a qu\u00edmica experimental. S\u00e3o Paulo: McGraw-Hill, 1990.\"\n$footer = \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n\n$markerIndex = -1\n$footerIndex = -1\n$n = $d.Paragraphs.Count\nfor ($i = 1; $i -le $n; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd(\"`r`a\")\n    if ($t -eq $marker) {\n        $markerIndex = $i\n    }\n    if ($t -eq $footer) {\n        $footerIndex = $i\n        break\n    }\n}\n\nif ($markerIndex -eq -1 -or $footerIndex -eq -1 -or $footerIndex -le $markerIndex) {\n    throw \"Could not locate the paragraphs to delete.\"\n}\n\n# Delete from the bottom up so earlier indices stay valid.\nfor ($i = $footerIndex; $i -gt $markerIndex; $i--) {\n    $d.Paragraphs.Item($i).Range.Delete()\n}\n"}
